# Update cryptocurrency market data scraped on Sun Jul 30 05:38:24 UTC 2023
# (GitHub Actions "cryptos list" refresh).
#
# Column D ("Price") and column E ("Volume(1h)") are refreshed with the
# latest snapshot values for every listed coin. Rows 45/46 additionally
# swap which coin (PaxDollar / BabyDogeCoin) occupies which rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source feed always stores Price/Volume as text (e.g. "29.362.06",
# "242.83", "  -0.02%  "), even when a value happens to look like a plain
# number. Force the Text number format before writing such values so Excel
# does not silently reinterpret them as numeric/date types.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$ws.Range("D2").Value = '29.362.06'
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").Value = '1.879.66'
$ws.Range("E3").Value = '  +0.28%  '

Set-TextValue $ws.Range("D5") '0.7122'
$ws.Range("E5").Value = '  -0.01%  '

Set-TextValue $ws.Range("D6") '242.83'

Set-TextValue $ws.Range("D8") '0.08041'
$ws.Range("E8").Value = '  +3.04%  '

Set-TextValue $ws.Range("D9") '0.3179'
$ws.Range("E9").Value = '  +2.04%  '

Set-TextValue $ws.Range("D10") '25.11'
$ws.Range("E10").Value = '  -0.31%  '

Set-TextValue $ws.Range("D11") '0.08341'
$ws.Range("E11").Value = '  -1.25%  '

$ws.Range("D12").Value = '1.896.59'
$ws.Range("E12").Value = '  +1.16%  '

Set-TextValue $ws.Range("D13") '5.264'
$ws.Range("E13").Value = '  +0.49%  '

Set-TextValue $ws.Range("D14") '94.89'
$ws.Range("E14").Value = '  +4.16%  '

Set-TextValue $ws.Range("D15") '0.7182'
$ws.Range("E15").Value = '  +0.67%  '

Set-TextValue $ws.Range("D16") '6.362'
$ws.Range("E16").Value = '  +4.88%  '

Set-TextValue $ws.Range("D17") '0.000008631'
$ws.Range("E17").Value = '  +4.86%  '

$ws.Range("D18").Value = '29.376.96'
$ws.Range("E18").Value = '  +0.01%  '

Set-TextValue $ws.Range("D19") '243.29'
$ws.Range("E19").Value = '  +0.99%  '

$ws.Range("D20").Value = '2.149.02'
$ws.Range("E20").Value = '  +1.26%  '

Set-TextValue $ws.Range("D21") '13.35'
$ws.Range("E21").Value = '  +0.78%  '

Set-TextValue $ws.Range("D22") '1.001'
$ws.Range("E22").Value = '  +0.19%  '

Set-TextValue $ws.Range("D23") '7.835'
$ws.Range("E23").Value = '  +0.57%  '

$ws.Range("E24").Value = '  +0.10%  '

Set-TextValue $ws.Range("D25") '0.1575'
$ws.Range("E25").Value = '  -1.18%  '

Set-TextValue $ws.Range("D26") '9.105'
$ws.Range("E26").Value = '  +0.38%  '

Set-TextValue $ws.Range("D27") '163.34'
$ws.Range("E27").Value = '  +0.06%  '

Set-TextValue $ws.Range("D28") '18.62'
$ws.Range("E28").Value = '  +0.30%  '

Set-TextValue $ws.Range("D29") '1.509'
$ws.Range("E29").Value = '  -0.16%  '

Set-TextValue $ws.Range("D30") '4.439'
$ws.Range("E30").Value = '  +0.38%  '

Set-TextValue $ws.Range("D31") '4.353'
$ws.Range("E31").Value = '  +0.59%  '

Set-TextValue $ws.Range("D32") '1.205'
$ws.Range("E32").Value = '  -6.67%  '

Set-TextValue $ws.Range("D33") '0.05413'
$ws.Range("E33").Value = '  +1.98%  '

Set-TextValue $ws.Range("D34") '1.948'
$ws.Range("E34").Value = '  +0.46%  '

Set-TextValue $ws.Range("D35") '0.7731'
$ws.Range("E35").Value = '  +3.78%  '

Set-TextValue $ws.Range("D36") '1.186'
$ws.Range("E36").Value = '  +0.61%  '

Set-TextValue $ws.Range("D37") '2.692'
$ws.Range("E37").Value = '  -0.19%  '

Set-TextValue $ws.Range("D38") '0.01892'
$ws.Range("E38").Value = '  +1.27%  '

$ws.Range("D39").Value = '1.269.62'
$ws.Range("E39").Value = '  +3.48%  '

Set-TextValue $ws.Range("D40") '2.753'
$ws.Range("E40").Value = '  +0.97%  '

Set-TextValue $ws.Range("D41") '6.506'
$ws.Range("E41").Value = '  -0.01%  '

Set-TextValue $ws.Range("D42") '114.01'
$ws.Range("E42").Value = '  +2.92%  '

Set-TextValue $ws.Range("D43") '0.9109'
$ws.Range("E43").Value = '  +2.16%  '

Set-TextValue $ws.Range("D44") '74.70'
$ws.Range("E44").Value = '  +2.49%  '

$ws.Range("D47").Value = '2.037.17'
$ws.Range("E47").Value = '  +0.85%  '

Set-TextValue $ws.Range("D48") '1.809'
$ws.Range("E48").Value = '  -0.13%  '

Set-TextValue $ws.Range("D49") '0.5224'
$ws.Range("E49").Value = '  +0.21%  '

Set-TextValue $ws.Range("D50") '9.559'
$ws.Range("E50").Value = '  +1.27%  '

Set-TextValue $ws.Range("D51") '0.4383'
$ws.Range("E51").Value = '  +1.36%  '

# Rows 45/46: the two coins swap positions (BabyDogeCoin now ranks above
# PaxDollar) and both get refreshed Price/Volume figures.
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D45") "0.00000000131"
$ws.Range("E45").Value = "  +6.52%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D46") "1.001"
$ws.Range("E46").Value = "  +0.13%  "
